# Update row 2 (Raymond ad text) - refresh ad copy and append a new amenities line.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Sponsored`nPresenting 2&3 BHK in Thane W - 40+ Indoor & Outdoor Amenities`nraymondtenxera.com`nhttps://www.raymondtenxera.com › thane › project`nExperience a futuristic lifestyle with Raymond Realty's Spacious 2 BHK homes in Thane West. Prime Connectivity: School - 5 Mins | Metro - 3 Mins | Hospital - 2 Mins | Mall - 2 Mins.`nLocation Map · Configuration · Overview · View Amenities · Location Advantages · About Us"

# Update row 3 (Runwal ad text) - swap last line of call-to-action links.
$ws.Range("A3").Value = "Sponsored`n2 BHK Projects in Thane West | Starts at ₹93 Lacs* by Runwal`nlandsend.runwal.com`nhttp://landsend.runwal.com › projects › thane`nTake Advantage of the Umbrella Offer: 2 BHK Flats Starts at ₹93L* at Lands End by Runwal. Book Your Dream Home at Runwal Lands End And Avail Pay 10% Now & 90% On Possession Plan. Kendriya Vidyalaya :2Mins. Amenities: Mini Theatre, Putting Golf.`nView Location · View Gallery · Lands End by Runwal · Project Configuration"

# Row 4 is now a Lodha Group ad (was Godrej); advertiser also changes.
$ws.Range("A4").Value = "Sponsored`nWorld-class 2 BHKs in Thane | 1,2,3 BHK by Lodha® in Thane`nLodha Group`nhttps://www.lodhagroup.in`nWorld-class 1 BHKs in Thane by India's #1 real estate developer. Building a better life. World-class homes in Thane by India's #1 real estate developer. Building a better life.`nPrice · Location · About · Amenities"
$ws.Range("B4").Value = "Macrotech Developers Limited"

# Insert a new row before the old row 5, pushing the Dosti row down to row 6,
# and populate it with the new Godrej Ascend ad. The previous row 4 advertiser
# (Madison Communications) now lands on this inserted row.
$ws.Rows("5").Insert()

$ws.Range("A5").Value = "Sponsored`nGodrej Ascend, Kolshet, Thane - 2&3 BHK at ₹1.09Cr+*(All Incl)`nGodrej Properties`nhttps://www.godrejproperties.com`nLive a grand life at Godrej Ascend, Thane | 2&3 BHK at ₹1.09Cr+(All Incl)*"
$ws.Range("B5").Value = "Madison Communications Private Limited"
$ws.Range("C5").Value = "India"

# Row 6 (formerly row 5) - Dosti Willow ad text rewritten.
$ws.Range("A6").Value = "Sponsored`n2 & 3 BHK Homes at Thane`ndostiwillow.com`nhttps://www.dostiwillow.com › homes › thane`nAt Balkum, Thane (W) Dosti Willow at Thane | Project by Dosti™ Realty | Thane's Signature Living Expanded. Lifestyle | Education | Sport | Spirituality | Recreation | Healthcare | Connectivity."
$ws.Range("B6").Value = "DOSTI REALTY LIMITED"
$ws.Range("C6").Value = "India"
